$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new DAMSLTag (col I) / DialogAct (col J) values
$updates = @{
    2   = @("sd", "Statement-non-opinion")
    51  = @("sd", "Statement-non-opinion")
    69  = @("sd", "Statement-non-opinion")
    101 = @("sd", "Statement-non-opinion")
    127 = @("sd", "Statement-non-opinion")
    128 = @("sd", "Statement-non-opinion")
    132 = @("sv", "Statement-opinion")
    134 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
